$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update numeric values (map/stage numbers adjusted)
$ws.Range("D6").Value = 130
$ws.Range("E6").Value = 120
$ws.Range("D8").Value = 130
$ws.Range("E8").Value = 120

# Update the active selection on the sheet from F13 to F12
$ws.Range("F12").Select()
